# E suite.xlsx - "taking latest changes and appending new profile scripts"
#
# This script updates the "Test Cases" worksheet:
#  - Row 28 (TestCase_E27) gets an extra Jira id appended and its description
#    is extended with a third clause.
#  - Nine new watchlist/profile test cases are appended as rows 31-35
#    (TestCase_E30 .. TestCase_E34).
#  - The worksheet dimension grows from A1:E30 to A1:E35.
#  - The active selection / scroll position of the sheet is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Update TestCase_E27 (row 28): new Jira id + extended description,
#    and the row grows taller to fit the extra wrapped text (30 -> 45).
# ---------------------------------------------------------------------
$ws.Range("B28").Value = "OPQA-314 ||OPQA-317`n||OPQA-327"
$ws.Range("C28").Value = "Verify that user is able to name the watchlists||Verify that a user can add description to his watchlist||Verify that watchlist name is customizable"
$ws.Rows.Item(28).RowHeight = 45

# ---------------------------------------------------------------------
# 2. Append the new test cases as rows 31-35. Column C alternates between
#    two existing "wrap text" cell styles already used in the sheet, so
#    clone the formatting from the nearest existing row of each kind
#    (row 29 uses style s="4", row 30 uses style s="3") instead of
#    re-building the style by hand. Columns A, B, D and E reuse the plain
#    style already used higher up in the table (e.g. row 2).
# ---------------------------------------------------------------------
function Clone-CellFormat($srcAddress, $dstAddress) {
    $ws.Range($srcAddress).Copy() | Out-Null
    $ws.Range($dstAddress).PasteSpecial(-4122) | Out-Null # xlPasteFormats
}

# Row 31
Clone-CellFormat "A2" "A31"
Clone-CellFormat "B2" "B31"
Clone-CellFormat "C29" "C31"
Clone-CellFormat "D2" "D31"
Clone-CellFormat "E2" "E31"

# Row 32
Clone-CellFormat "A2" "A32"
Clone-CellFormat "B2" "B32"
Clone-CellFormat "C30" "C32"
Clone-CellFormat "D2" "D32"
Clone-CellFormat "E2" "E32"

# Row 33
Clone-CellFormat "A2" "A33"
Clone-CellFormat "B2" "B33"
Clone-CellFormat "C29" "C33"
Clone-CellFormat "D2" "D33"
Clone-CellFormat "E2" "E33"

# Row 34
Clone-CellFormat "A2" "A34"
Clone-CellFormat "B2" "B34"
Clone-CellFormat "C30" "C34"
Clone-CellFormat "D2" "D34"
Clone-CellFormat "E2" "E34"

# Row 35
Clone-CellFormat "A2" "A35"
Clone-CellFormat "B2" "B35"
Clone-CellFormat "C30" "C35"
Clone-CellFormat "D2" "D35"
Clone-CellFormat "E2" "E35"

$excel.CutCopyMode = $false

# Values for the newly appended rows.
$ws.Range("A31").Value = "TestCase_E30"
$ws.Range("B31").Value = "OPQA-324"
$ws.Range("C31").Value = "Verify that a user has 1 watchlist by default once we try to watch an item"
$ws.Range("D31").Value = "Y"
$ws.Range("E31").Value = "PASS"

$ws.Range("A32").Value = "TestCase_E31"
$ws.Range("B32").Value = "OPQA-326"
$ws.Range("C32").Value = "Verify that user is able to have a watchlist with 0 item under it"
$ws.Range("D32").Value = "Y"
$ws.Range("E32").Value = "PASS"

$ws.Range("A33").Value = "TestCase_E32"
$ws.Range("B33").Value = "OPQA-328"
$ws.Range("C33").Value = "Verify that every user watchlist is private by default"
$ws.Range("D33").Value = "Y"
$ws.Range("E33").Value = "PASS"

$ws.Range("A34").Value = "TestCase_E33"
$ws.Range("B34").Value = "OPQA-321"
$ws.Range("C34").Value = "Verify that anyone can see the public watchlists of a user on user's profile page"
$ws.Range("D34").Value = "Y"
$ws.Range("E34").Value = "PASS"

$ws.Range("A35").Value = "TestCase_E34"
$ws.Range("B35").Value = "OPQA-329"
$ws.Range("C35").Value = "Verify that no one can see the private watchlists of a user on user's profile page"
$ws.Range("D35").Value = "Y"
$ws.Range("E35").Value = "PASS"

# ---------------------------------------------------------------------
# 3. Update the view: scroll so column C is the left-most visible
#    column, and select E2:E35 (active cell E2).
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E2:E35").Select()

# Best-effort: shrink the workbook window height recorded for the book.
$excel.ActiveWindow.Height = 8460
